$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rngA1 = $wsHoja1.Range("A1")
$oldText = $rngA1.Value2
$newText = $oldText.Replace(
    "1000 Bs = 10.22 = 42660.02 pesos",
    "1000 Bs = 10.24 = 42374.62 pesos"
).Replace(
    "42660.02 pesos = 10.21 = 959.34 Bs",
    "42374.62 pesos = 10.18 = 959.13 Bs"
)
$rngA1.Value = $newText

# --- Sheet "tasas": update the rate table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 97.7
$wsTasas.Range("O10").Value = 4140
$wsTasas.Range("N12").Value = 4162
$wsTasas.Range("O12").Value = 94.205
